$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header columns G:J mirror B:E headers (period = 1 / 5 / 10 / 30)
$ws.Range("G2").Value = "period = 1"
$ws.Range("H2").Value = "period = 5"
$ws.Range("I2").Value = "period = 10"
$ws.Range("J2").Value = "period = 30"

# New data block G3:J7
$ws.Range("G3").Value = 65.47
$ws.Range("H3").Value = 55.85
$ws.Range("I3").Value = 54.78
$ws.Range("J3").Value = 36.47

$ws.Range("G4").Value = 56.98
$ws.Range("H4").Value = 51.35
$ws.Range("I4").Value = 52.72
$ws.Range("J4").Value = 33.33

$ws.Range("G5").Value = 65.42
$ws.Range("H5").Value = 55.08
$ws.Range("I5").Value = 45.7
$ws.Range("J5").Value = 36.26

$ws.Range("G6").Value = 61.08
$ws.Range("H6").Value = 54.35
$ws.Range("I6").Value = 38.79
$ws.Range("J6").Value = 29.61

$ws.Range("G7").Value = 57.75
$ws.Range("H7").Value = 58
$ws.Range("I7").Value = 44.78
$ws.Range("J7").Value = 34.95

# New "Total" row with AVERAGE formulas
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Formula = "=AVERAGE(B3:B7)"
$ws.Range("C8:F8").Formula = "=AVERAGE(C3:C7)"
$ws.Range("F8").ClearContents()
$ws.Range("G8").Formula = "=AVERAGE(G3:G7)"
$ws.Range("H8").Formula = "=AVERAGE(H3:H7)"
$ws.Range("I8:J8").Formula = "=AVERAGE(I3:I7)"

# Column width adjustments for new columns H and I
$ws.Columns.Item(8).ColumnWidth = 11.43
$ws.Columns.Item(9).ColumnWidth = 13.92

# Update selection to match target (I11)
$ws.Range("I11").Select()
